# Adds a new date column ("25. 1. 2021") at the end of both data tables
# ("data" and "pocetR" sheets) and refreshes the "aktualizace" date in the
# title rows from 6. 1. 2022 to 1. 2. 2022.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": header row values are in columns A:AL (1..38); new data
# goes into column AM (39).
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

$newHeaderCol = 39   # AM
$lastHeaderCol = 38  # AL

# Header cell - copy format from the previous header cell (AL1) so the new
# header keeps the same bold/centered/bordered style, then set its value.
$wsData.Cells.Item(1, $lastHeaderCol).Copy() | Out-Null
$wsData.Cells.Item(1, $newHeaderCol).PasteSpecial(-4122) | Out-Null
$wsData.Cells.Item(1, $newHeaderCol).Value = "25. 1. 2021"

$dataValues = @{
    2  = 0.11
    3  = 0.12
    4  = 0.11
    5  = 0.12
    6  = 0.12
    7  = 0.12
    8  = 0.09
    9  = 0.1
    10 = 0.12
    11 = 0.08
    12 = 0.11
    13 = 0.12
    14 = 0.12
    15 = 0.17
    16 = 0.15
    17 = 0.1
    18 = 0.1
    19 = 0.1
    20 = 0.22
    21 = 0.14
    22 = 0.09
    23 = 0.2
    24 = 0.1
    25 = 0.1
    26 = 0.1
    27 = 0.07
    28 = 0.09
    29 = 0.18
}

foreach ($row in $dataValues.Keys) {
    $wsData.Cells.Item($row, $newHeaderCol).Value = $dataValues[$row]
}

# Update title on the final row (A30) with the new "aktualizace" date.
$wsData.Cells.Item(30, 1).Value = "Život během pandemie, Duševní zdraví, % respondentů celkově a ve skupinách, aktualizace 1. 2. 2022"

# ---------------------------------------------------------------------
# Sheet "pocetR": header row values are in columns A:AK (1..37); new data
# goes into column AL (38). (This sheet has no "ukazatel" column, so it is
# shifted one column to the left compared to "data".)
# ---------------------------------------------------------------------
$wsCount = $wb.Worksheets.Item("pocetR")

$newHeaderColCount = 38   # AL
$lastHeaderColCount = 37  # AK

$wsCount.Cells.Item(1, $lastHeaderColCount).Copy() | Out-Null
$wsCount.Cells.Item(1, $newHeaderColCount).PasteSpecial(-4122) | Out-Null
$wsCount.Cells.Item(1, $newHeaderColCount).Value = "25. 1. 2021"

$countValues = @{
    2  = 1815
    3  = 412
    4  = 1403
    5  = 304
    6  = 797
    7  = 109
    8  = 605
    9  = 880
    10 = 935
    11 = 220
    12 = 660
    13 = 293
    14 = 642
    15 = 161
    16 = 283
    17 = 356
    18 = 320
    19 = 695
    20 = 186
    21 = 324
    22 = 1305
    23 = 170
    24 = 630
    25 = 605
    26 = 302
    27 = 494
    28 = 759
    29 = 562
}

foreach ($row in $countValues.Keys) {
    $wsCount.Cells.Item($row, $newHeaderColCount).Value = $countValues[$row]
}

# Update title on the final row (A30) with the new "aktualizace" date.
$wsCount.Cells.Item(30, 1).Value = "Život během pandemie, Duševní zdraví, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 2. 2022"
